# Feria Lagunitas de Puerto Montt - Pomelo: insert a new weekly record
# at row 661, pushing the existing rows 661-714 down to 662-715.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 661 (shifts rows 661:714 -> 662:715,
# carries formatting down automatically, and grows the sheet dimension).
$ws.Rows("661:661").Insert()

# Populate the newly inserted row 661 with the new record's data.
$row = 661
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 45265
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100102
$ws.Cells.Item($row, 8).Value = "Cítricos"
$ws.Cells.Item($row, 9).Value = 100102006
$ws.Cells.Item($row, 10).Value = "Pomelo"
$ws.Cells.Item($row, 11).Value = "Start Ruby"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 100
$ws.Cells.Item($row, 14).Value = 15000
$ws.Cells.Item($row, 15).Value = 15000
$ws.Cells.Item($row, 16).Value = 15000
$ws.Cells.Item($row, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 1071
$ws.Cells.Item($row, 20).Value = 14
